$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style (bold, bordered, centered) used by the other header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# I0 / IF values for rows 2-29
$values = @(
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 9),
    @(1, 4),
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 6),
    @(1, 5),
    @(1, 3),
    @(2, 7),
    @(6, 7),
    @(5, 7),
    @(6, 8),
    @(4, 7),
    @(4, 6),
    @(4, 7),
    @(6, 7),
    @(2, 5),
    @(8, 8),
    @(6, 7),
    @(7, 8),
    @(5, 6),
    @(4, 5)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}

$ws.Range("A1:J29").Select()
